# Word COM-interop script reproducing the tracked edits for
# "skeleton based action recognition using few shot learning.docx"
#
# Summary of the real (visible) content changes made by the author:
#   1. Author byline runs re-typed/merged (no visible text change, but
#      Word re-flows the runs and drops the now-stale proofing marks).
#   2. Introduction paragraph reworded:
#        "One of the them"      -> "One of them"
#        "into a feature image."-> "into a skeleton-image."
#        new sentences inserted ("Different skeleton-image encoders ...
#        Inspired by the action movement mechanism,   ") right before
#        "The other one is ...", and "based on" -> "tend to",
#        "Recurrent" -> "recurrent".
#
# wdReplaceAll = 2 ; wdFindContinue (Wrap) = 1

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Byline: "Wangbin Ding1,2, Shengqin Lin2, Zhenze Dai2"
#    Re-key the plain-text runs so Word merges them the same way it
#    does whenever a user retypes/edits text inside an existing run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Wangbin Ding", $false, $false, $false, $false, $false, $true, 1, $false, "Wangbin Ding", 2)
$d.Content.Find.Execute(", Shen", $false, $false, $false, $false, $false, $true, 1, $false, ", Shen", 2)
$d.Content.Find.Execute("qin Lin", $false, $false, $false, $false, $false, $true, 1, $false, "qin Lin", 2)
$d.Content.Find.Execute(", Zhenze Dai", $false, $false, $false, $false, $false, $true, 1, $false, ", Zhenze Dai", 2)

# ---------------------------------------------------------------------
# 2. Introduction paragraph rewording.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("One of the them", $false, $false, $false, $false, $false, $true, 1, $false, "One of them", 2)

$d.Content.Find.Execute("into a feature image.", $false, $false, $false, $false, $false, $true, 1, $false, "into a skeleton-image.", 2)

$d.Content.Find.Execute(
    "The other one is based on the RNN network, whose Recurrent structure",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Different skeleton-image encoders are proposed to capture the feature of actions. Inspired by the action movement mechanism,   The other one is tend to the RNN network, whose recurrent structure",
    2)
